$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 508 (this shifts all rows from 508 down
# to 510 onward, preserving formatting/styles from the row above as Excel
# normally does on row insert).
$ws.Rows("508:509").Insert()

# Populate the two newly inserted rows with the new weekly price report.

# Row 508 - Coliflor, Primera
$ws.Range("A508").Value = 8
$ws.Range("B508").Value = "Terminal La Palmera de La Serena"
$ws.Range("C508").Value = "Coquimbo"
$ws.Range("D508").Value = 44585
$ws.Range("E508").Value = 4
$ws.Range("F508").Value = 100112008
$ws.Range("G508").Value = "Coliflor"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 2400
$ws.Range("K508").Value = 750
$ws.Range("L508").Value = 800
$ws.Range("M508").Value = 775
$ws.Range("N508").Value = "`$/unidad"
$ws.Range("O508").Value = "Provincia del Elquí"
$ws.Range("P508").Value = 775
$ws.Range("Q508").Value = 1
$ws.Range("R508").Value = "Hortaliza"

# Row 509 - Coliflor, Segunda
$ws.Range("A509").Value = 8
$ws.Range("B509").Value = "Terminal La Palmera de La Serena"
$ws.Range("C509").Value = "Coquimbo"
$ws.Range("D509").Value = 44585
$ws.Range("E509").Value = 4
$ws.Range("F509").Value = 100112008
$ws.Range("G509").Value = "Coliflor"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Segunda"
$ws.Range("J509").Value = 1700
$ws.Range("K509").Value = 650
$ws.Range("L509").Value = 700
$ws.Range("M509").Value = 675
$ws.Range("N509").Value = "`$/unidad"
$ws.Range("O509").Value = "Provincia del Elquí"
$ws.Range("P509").Value = 675
$ws.Range("Q509").Value = 1
$ws.Range("R509").Value = "Hortaliza"

# Ensure the date cells keep the date number format used elsewhere in column D.
$ws.Range("D508:D509").NumberFormat = $ws.Range("D510").NumberFormat()
